# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" sheet and the "全部类型" sheet (they carry duplicate data).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 11591
    $ws.Range("F3").Value = 11133
    $ws.Range("F8").Value = 67
    $ws.Range("F11").Value = 10674
    $ws.Range("F21").Value = 10881
}
